$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(66539.3665, 45000, 25000, 20000, 10000, 10000, 0, 16022.73300000001, 1),
    @(63920.688, 45000, 25000, 20000, 10000, 10000, 0, 13533.376, 2),
    @(60272.176, 22500, 12500, 10000, 5000, 5000, 0, 32432.352, 3),
    @(58962.181, 22500, 12500, 10000, 5000, 5000, 0, 31155.36199999999, 4),
    @(59797.17, 32500, 12500, 20000, 10000, 10000, 0, 21958.34, 5),
    @(65642.2855, 22500, 12500, 10000, 5000, 5000, 0, 37719.571, 6),
    @(64173.893, 22500, 12500, 10000, 5000, 5000, 0, 35798.78599999999, 7),
    @(73269.5725, 22500, 12500, 10000, 5000, 5000, 0, 43555.14499999999, 8),
    @(86257.982, 65000, 25000, 20000, 10000, 10000, 20000, 12790.96400000001, 9),
    @(90482.3835, 65000, 25000, 20000, 10000, 10000, 20000, 11841.76699999999, 10),
    @(94548.98, 65000, 25000, 20000, 10000, 10000, 20000, 13785.95999999999, 11),
    @(96248.9605, 65000, 25000, 20000, 10000, 10000, 20000, 16176.921, 12),
    @(93999.3325, 65000, 25000, 20000, 10000, 10000, 20000, 13739.66500000001, 13),
    @(97613.256, 65000, 25000, 20000, 10000, 10000, 20000, 16797.51199999999, 14),
    @(103207.6725, 65000, 25000, 20000, 10000, 10000, 20000, 22252.345, 15),
    @(53761.9265, 65000, 25000, 20000, 10000, 10000, 20000, -27352.147, 16),
    @(86334.611, 65000, 25000, 20000, 10000, 10000, 20000, 4950.222000000009, 17),
    @(82324.363, 65000, 25000, 20000, 10000, 10000, 20000, 1391.725999999995, 18),
    @(81334.0935, 65000, 25000, 20000, 10000, 10000, 20000, 1297.187000000005, 19),
    @(78578.4975, 65000, 25000, 20000, 10000, 10000, 20000, 247.9949999999953, 20),
    @(77749.1885, 65000, 25000, 20000, 10000, 10000, 20000, 978.3770000000077, 21),
    @(74837.9245, 65000, 25000, 20000, 10000, 10000, 20000, 525.8489999999947, 22),
    @(54075.2015, 65000, 25000, 20000, 10000, 10000, 20000, -17698.59699999999, 23),
    @(14561.85, 65000, 25000, 20000, 10000, 10000, 20000, -55980.281, 24)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}
